$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dateCreated origin-info fragment's encoding attribute from
# iso8601 to w3cdtf.
$ws.Range("S1").Value = '<mods:originInfo><mods:dateCreated encoding="w3cdtf">'

# Reflect the author's last on-screen selection/scroll position.
$ws.Range("S1").Select()
